$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DATE column values ("2.1.19") look like an ambiguous d.m.yy date, so force
# the cells to Text first -- matches existing rows like "16.12.18", which are
# plain shared-string text, not real date serials.
$ws.Range("B10:B13").NumberFormat = "@"

# Append four new commit-log rows (rows 10-13)
$ws.Range("A10").Value = "Omri"
$ws.Range("B10").Value = "2.1.19"
$ws.Range("C10").Value = "vf2pf.v"
$ws.Range("D10").Value = "inserting basic translation for spesific case in order to check functionality"

$ws.Range("A11").Value = "Omri"
$ws.Range("B11").Value = "2.1.19"
$ws.Range("C11").Value = "address_translation.v"
$ws.Range("D11").Value = "selection between regular pf , vf2pf , 4k translation added , shourtcuted to the output"

$ws.Range("A12").Value = "Ori"
$ws.Range("B12").Value = "2.1.19"
$ws.Range("C12").Value = "spi.v"
$ws.Range("D12").Value = "connecting cache req to nvm read stage in spi state machine "

$ws.Range("A13").Value = "Ori"
$ws.Range("B13").Value = "2.1.19"
$ws.Range("C13").Value = "spi.v"
$ws.Range("D13").Value = "collect read data from nvm interface "

# The DATE column never visually used the Text number-format (the source
# rows are plain "General" cells that merely contain date-like strings), so
# drop back to the Normal cell style now that the text entry is locked in.
$ws.Range("B10:B13").Style = "Normal"

# Select the new last cell, mirroring the recorded user action
$ws.Range("D13").Select()

# Autofit columns A:D to match the recorded column widths
$ws.Columns("A:D").AutoFit() | Out-Null
